$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 12,20
$arr[0,0] = "ECs"
$arr[0,1] = "Lta"
$arr[0,2] = "Tnfrsf1a"
$arr[0,3] = "ECs"
$arr[0,4] = 1
$arr[0,5] = 0.3333333333333333
$arr[0,6] = 0.103879
$arr[0,7] = 0.311637
$arr[0,8] = 0.1899090787212519
$arr[0,9] = 0.1899090787212519
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 38.45264233333334
$arr[0,13] = 115.357927
$arr[0,14] = 0.2975040117664333
$arr[0,15] = 0.2975040117664332
$arr[0,16] = 3.994422032944334
$arr[0,17] = 35.949798296499
$arr[0,18] = 0.05649871279043983
$arr[0,19] = 0.05649871279043982
$arr[1,0] = "ECs"
$arr[1,1] = "Lta"
$arr[1,2] = "Tnfrsf1a"
$arr[1,3] = "FAPs"
$arr[1,4] = 1
$arr[1,5] = 0.3333333333333333
$arr[1,6] = 0.103879
$arr[1,7] = 0.311637
$arr[1,8] = 0.1899090787212519
$arr[1,9] = 0.1899090787212519
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 47.750315
$arr[1,13] = 143.250945
$arr[1,14] = 0.3694391181876273
$arr[1,15] = 0.3694391181876272
$arr[1,16] = 4.960254971885
$arr[1,17] = 44.642294746965
$arr[1,18] = 0.070159842578604
$arr[1,19] = 0.07015984257860398
$arr[2,0] = "ECs"
$arr[2,1] = "Lta"
$arr[2,2] = "Tnfrsf1a"
$arr[2,3] = "MuSCs"
$arr[2,4] = 1
$arr[2,5] = 0.3333333333333333
$arr[2,6] = 0.103879
$arr[2,7] = 0.311637
$arr[2,8] = 0.1899090787212519
$arr[2,9] = 0.1899090787212519
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 18.63107466666667
$arr[2,13] = 55.893224
$arr[2,14] = 0.1441466469015163
$arr[2,15] = 0.1441466469015162
$arr[2,16] = 1.935377405298667
$arr[2,17] = 17.418396647688
$arr[2,18] = 0.02737475691382456
$arr[2,19] = 0.02737475691382455
$arr[3,0] = "ECs"
$arr[3,1] = "Lta"
$arr[3,2] = "Tnfrsf1a"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 1
$arr[3,5] = 0.3333333333333333
$arr[3,6] = 0.103879
$arr[3,7] = 0.311637
$arr[3,8] = 0.1899090787212519
$arr[3,9] = 0.1899090787212519
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 24.41680433333333
$arr[3,13] = 73.25041299999999
$arr[3,14] = 0.1889102231444233
$arr[3,15] = 0.1889102231444233
$arr[3,16] = 2.536393217342333
$arr[3,17] = 22.827538956081
$arr[3,18] = 0.03587576643838355
$arr[3,19] = 0.03587576643838355
$arr[4,0] = "FAPs"
$arr[4,1] = "Lta"
$arr[4,2] = "Tnfrsf1a"
$arr[4,3] = "ECs"
$arr[4,4] = 2
$arr[4,5] = 0.6666666666666666
$arr[4,6] = 0.257314
$arr[4,7] = 0.771942
$arr[4,8] = 0.4704152396738534
$arr[4,9] = 0.4704152396738534
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 38.45264233333334
$arr[4,13] = 115.357927
$arr[4,14] = 0.2975040117664333
$arr[4,15] = 0.2975040117664332
$arr[4,16] = 9.894403209359334
$arr[4,17] = 89.04962888423401
$arr[4,18] = 0.1399504209990396
$arr[4,19] = 0.1399504209990396
$arr[5,0] = "FAPs"
$arr[5,1] = "Lta"
$arr[5,2] = "Tnfrsf1a"
$arr[5,3] = "FAPs"
$arr[5,4] = 2
$arr[5,5] = 0.6666666666666666
$arr[5,6] = 0.257314
$arr[5,7] = 0.771942
$arr[5,8] = 0.4704152396738534
$arr[5,9] = 0.4704152396738534
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 47.750315
$arr[5,13] = 143.250945
$arr[5,14] = 0.3694391181876273
$arr[5,15] = 0.3694391181876272
$arr[5,16] = 12.28682455391
$arr[5,17] = 110.58142098519
$arr[5,18] = 0.1737897913271297
$arr[5,19] = 0.1737897913271297
$arr[6,0] = "FAPs"
$arr[6,1] = "Lta"
$arr[6,2] = "Tnfrsf1a"
$arr[6,3] = "MuSCs"
$arr[6,4] = 2
$arr[6,5] = 0.6666666666666666
$arr[6,6] = 0.257314
$arr[6,7] = 0.771942
$arr[6,8] = 0.4704152396738534
$arr[6,9] = 0.4704152396738534
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 18.63107466666667
$arr[6,13] = 55.893224
$arr[6,14] = 0.1441466469015163
$arr[6,15] = 0.1441466469015162
$arr[6,16] = 4.794036346778666
$arr[6,17] = 43.14632712100801
$arr[6,18] = 0.0678087794503591
$arr[6,19] = 0.06780877945035908
$arr[7,0] = "FAPs"
$arr[7,1] = "Lta"
$arr[7,2] = "Tnfrsf1a"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.257314
$arr[7,7] = 0.771942
$arr[7,8] = 0.4704152396738534
$arr[7,9] = 0.4704152396738534
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 24.41680433333333
$arr[7,13] = 73.25041299999999
$arr[7,14] = 0.1889102231444233
$arr[7,15] = 0.1889102231444233
$arr[7,16] = 6.282785590227332
$arr[7,17] = 56.545070312046
$arr[7,18] = 0.08886624789732502
$arr[7,19] = 0.088866247897325
$arr[8,0] = "Resolving-Mac"
$arr[8,1] = "Lta"
$arr[8,2] = "Tnfrsf1a"
$arr[8,3] = "ECs"
$arr[8,4] = 2
$arr[8,5] = 0.6666666666666666
$arr[8,6] = 0.1858003333333333
$arr[8,7] = 0.557401
$arr[8,8] = 0.3396756816048946
$arr[8,9] = 0.3396756816048946
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 38.45264233333334
$arr[8,13] = 115.357927
$arr[8,14] = 0.2975040117664333
$arr[8,15] = 0.2975040117664332
$arr[8,16] = 7.144513763080779
$arr[8,17] = 64.30062386772701
$arr[8,18] = 0.1010548779769538
$arr[8,19] = 0.1010548779769538
$arr[9,0] = "Resolving-Mac"
$arr[9,1] = "Lta"
$arr[9,2] = "Tnfrsf1a"
$arr[9,3] = "FAPs"
$arr[9,4] = 2
$arr[9,5] = 0.6666666666666666
$arr[9,6] = 0.1858003333333333
$arr[9,7] = 0.557401
$arr[9,8] = 0.3396756816048946
$arr[9,9] = 0.3396756816048946
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 47.750315
$arr[9,13] = 143.250945
$arr[9,14] = 0.3694391181876273
$arr[9,15] = 0.3694391181876272
$arr[9,16] = 8.872024443771668
$arr[9,17] = 79.84821999394501
$arr[9,18] = 0.1254894842818935
$arr[9,19] = 0.1254894842818935
$arr[10,0] = "Resolving-Mac"
$arr[10,1] = "Lta"
$arr[10,2] = "Tnfrsf1a"
$arr[10,3] = "MuSCs"
$arr[10,4] = 2
$arr[10,5] = 0.6666666666666666
$arr[10,6] = 0.1858003333333333
$arr[10,7] = 0.557401
$arr[10,8] = 0.3396756816048946
$arr[10,9] = 0.3396756816048946
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 18.63107466666667
$arr[10,13] = 55.893224
$arr[10,14] = 0.1441466469015163
$arr[10,15] = 0.1441466469015162
$arr[10,16] = 3.461659883424889
$arr[10,17] = 31.154938950824
$arr[10,18] = 0.04896311053733261
$arr[10,19] = 0.04896311053733261
$arr[11,0] = "Resolving-Mac"
$arr[11,1] = "Lta"
$arr[11,2] = "Tnfrsf1a"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 2
$arr[11,5] = 0.6666666666666666
$arr[11,6] = 0.1858003333333333
$arr[11,7] = 0.557401
$arr[11,8] = 0.3396756816048946
$arr[11,9] = 0.3396756816048946
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 24.41680433333333
$arr[11,13] = 73.25041299999999
$arr[11,14] = 0.1889102231444233
$arr[11,15] = 0.1889102231444233
$arr[11,16] = 4.536650384068111
$arr[11,17] = 40.829853456613
$arr[11,18] = 0.06416820880871472
$arr[11,19] = 0.06416820880871471

$rng = $ws.Range("A2:T13")
$rng.Value2 = $arr
